$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)

# ---------------------------------------------------------------------
# 1) AutoShape 5 (id=10246) -- connector above "Text Box 6": shift/stretch
# ---------------------------------------------------------------------
$cxn5 = $s.Shapes.Item(5)
$cxn5.Top = 333.70762645669294
$cxn5.Height = 35.79165460629921

# ---------------------------------------------------------------------
# 2) Text Box 6 (id=10247): reposition/resize and shrink font to 16pt
# ---------------------------------------------------------------------
$tb6 = $s.Shapes.Item(6)
$tb6.Left = 47.83047303149606
$tb6.Width = 624.2875672047244
$tb6.Height = 26.70834645669291
$tb6.TextFrame.TextRange.Font.Size = 16

# ---------------------------------------------------------------------
# 3) Remove the old "Group 6" (id=7) holding the token text boxes
# ---------------------------------------------------------------------
$s.Shapes.Item("Group 6").Delete()

# ---------------------------------------------------------------------
# 4) Add the replacement group of token text boxes (smaller font,
#    merged runs, new wording/positions) at the bottom of the slide
# ---------------------------------------------------------------------
$t1 = $s.Shapes.AddTextbox(1, 24.0, 420.0, 159.62503937007875, 26.75)
$t1.TextFrame.TextRange.Text = "identifier [" + [char]8220 + "y" + [char]8221 + ", (1, 1)]"
$t1.TextFrame.TextRange.Font.Size = 16

$t2 = $s.Shapes.AddTextbox(1, 188.7812598425197, 420.0, 80.12503937007874, 26.75)
$t2.TextFrame.TextRange.Text = ":= [(1, 3)]"
$t2.TextFrame.TextRange.Font.Size = 16

$t3 = $s.Shapes.AddTextbox(1, 274.0625196850394, 420.0, 159.62503937007875, 26.75)
$t3.TextFrame.TextRange.Text = "identifier [" + [char]8220 + "x" + [char]8221 + ", (1, 6)]"
$t3.TextFrame.TextRange.Font.Size = 16

$t4 = $s.Shapes.AddTextbox(1, 438.84377952755904, 420.0, 75.62503937007874, 26.75)
$t4.TextFrame.TextRange.Text = "+ [(1, 8)]"
$t4.TextFrame.TextRange.Font.Size = 16

$t5 = $s.Shapes.AddTextbox(1, 519.6249606299212, 420.0, 170.37503937007875, 26.75)
$t5.TextFrame.TextRange.Text = "intLiteral [(" + [char]8220 + "1" + [char]8221 + ", (1, 10)]"
$t5.TextFrame.TextRange.Font.Size = 16

$range = $s.Shapes.Range(@($t1.Name, $t2.Name, $t3.Name, $t4.Name, $t5.Name))
$grp = $range.Group()
$grp.Left = 27.0
$grp.Top = 450.58323673228347
